# edit.ps1
# Applies 'want to go' count (column F) updates across the
# 展览(1) / 演出(2) / 全部类型(4) worksheets, per the source diff.
# 本地生活(3) sheet has no changes.

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F3").Value = 1170
$ws.Range("F4").Value = 1599
$ws.Range("F5").Value = 179
$ws.Range("F6").Value = 179
$ws.Range("F7").Value = 31
$ws.Range("F8").Value = 1554
$ws.Range("F9").Value = 3155
$ws.Range("F10").Value = 698
$ws.Range("F11").Value = 1878
$ws.Range("F12").Value = 1843
$ws.Range("F13").Value = 908
$ws.Range("F14").Value = 312
$ws.Range("F16").Value = 1529
$ws.Range("F19").Value = 38
$ws.Range("F20").Value = 1319
$ws.Range("F21").Value = 435
$ws.Range("F22").Value = 539
$ws.Range("F23").Value = 219
$ws.Range("F24").Value = 7837
$ws.Range("F25").Value = 7335
$ws.Range("F26").Value = 785
$ws.Range("F27").Value = 602
$ws.Range("F28").Value = 1740
$ws.Range("F29").Value = 104
$ws.Range("F30").Value = 273

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F5").Value = 33
$ws.Range("F6").Value = 107

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F5").Value = 1170
$ws.Range("F6").Value = 1599
$ws.Range("F7").Value = 179
$ws.Range("F8").Value = 179
$ws.Range("F10").Value = 31
$ws.Range("F11").Value = 1554
$ws.Range("F12").Value = 3155
$ws.Range("F13").Value = 698
$ws.Range("F14").Value = 1878
$ws.Range("F15").Value = 1843
$ws.Range("F16").Value = 908
$ws.Range("F17").Value = 312
$ws.Range("F19").Value = 1529
$ws.Range("F23").Value = 38
$ws.Range("F24").Value = 33
$ws.Range("F25").Value = 1319
$ws.Range("F26").Value = 435
$ws.Range("F27").Value = 539
$ws.Range("F28").Value = 219
$ws.Range("F29").Value = 7837
$ws.Range("F30").Value = 7336
$ws.Range("F31").Value = 785
$ws.Range("F32").Value = 602
$ws.Range("F33").Value = 1740
$ws.Range("F34").Value = 107
$ws.Range("F36").Value = 104
$ws.Range("F37").Value = 273

